# Auto-generated COM-interop script applying the scraped schedule update
# (Línea 141 - commit: 'Horarios actualizados Línea 141 - 253')
$wb = $excel.ActiveWorkbook

### Sheet: LP1912
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 04:54:03"
$ws.Range("A3").Value = "Total filas: 26"
$ws.Range("A6").Value = "04:18:53"
$ws.Range("B6").Value = "04:45"
$ws.Range("C6").Value = "215A_EL PATO"
$ws.Range("D6").Value = 27
$ws.Range("E6").Value = "LP1912"
$ws.Range("A7").Value = "04:40:33"
$ws.Range("B7").Value = "04:46"
$ws.Range("C7").Value = "215A_EL PATO"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = "LP1912"
$ws.Range("A8").Value = "04:18:53"
$ws.Range("B8").Value = "04:53"
$ws.Range("C8").Value = "11_ETCHEVERRY"
$ws.Range("D8").Value = 35
$ws.Range("E8").Value = "LP1912"
$ws.Range("A9").Value = "04:54:03"
$ws.Range("B9").Value = "04:55"
$ws.Range("C9").Value = "11_ETCHEVERRY"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "LP1912"
$ws.Range("A10").Value = "04:54:03"
$ws.Range("B10").Value = "05:15"
$ws.Range("C10").Value = "14_ABASTO"
$ws.Range("D10").Value = 21
$ws.Range("E10").Value = "LP1912"
$ws.Range("A11").Value = "04:18:53"
$ws.Range("B11").Value = "05:16"
$ws.Range("C11").Value = "17_ROMERO"
$ws.Range("D11").Value = 58
$ws.Range("E11").Value = "LP1912"
$ws.Range("A12").Value = "04:40:33"
$ws.Range("B12").Value = "05:16"
$ws.Range("C12").Value = "14_ABASTO"
$ws.Range("D12").Value = 36
$ws.Range("E12").Value = "LP1912"
$ws.Range("A13").Value = "04:18:53"
$ws.Range("B13").Value = "05:21"
$ws.Range("C13").Value = "23_HERNANDEZ"
$ws.Range("D13").Value = 63
$ws.Range("E13").Value = "LP1912"
$ws.Range("A14").Value = "04:40:33"
$ws.Range("B14").Value = "05:22"
$ws.Range("C14").Value = "23_HERNANDEZ"
$ws.Range("D14").Value = 42
$ws.Range("E14").Value = "LP1912"
$ws.Range("A15").Value = "04:18:53"
$ws.Range("B15").Value = "05:34"
$ws.Range("C15").Value = "215B_EL PATO"
$ws.Range("D15").Value = 76
$ws.Range("E15").Value = "LP1912"
$ws.Range("A16").Value = "04:54:03"
$ws.Range("B16").Value = "05:35"
$ws.Range("C16").Value = "215B_EL PATO"
$ws.Range("D16").Value = 41
$ws.Range("E16").Value = "LP1912"
$ws.Range("A17").Value = "04:18:53"
$ws.Range("B17").Value = "05:37"
$ws.Range("C17").Value = "14_ABASTO"
$ws.Range("D17").Value = 79
$ws.Range("E17").Value = "LP1912"
$ws.Range("A18").Value = "04:18:53"
$ws.Range("B18").Value = "05:46"
$ws.Range("C18").Value = "15_ABASTO"
$ws.Range("D18").Value = 88
$ws.Range("E18").Value = "LP1912"
$ws.Range("A19").Value = "04:40:33"
$ws.Range("B19").Value = "06:04"
$ws.Range("C19").Value = "16_SANTA ANA"
$ws.Range("D19").Value = 84
$ws.Range("E19").Value = "LP1912"
$ws.Range("A20").Value = "04:18:53"
$ws.Range("B20").Value = "06:07"
$ws.Range("C20").Value = "16_SANTA ANA"
$ws.Range("D20").Value = 109
$ws.Range("E20").Value = "LP1912"
$ws.Range("A21").Value = "04:18:53"
$ws.Range("B21").Value = "06:11"
$ws.Range("C21").Value = "215A_EL PATO"
$ws.Range("D21").Value = 113
$ws.Range("E21").Value = "LP1912"
$ws.Range("A22").Value = "04:54:03"
$ws.Range("B22").Value = "06:12"
$ws.Range("C22").Value = "215A_EL PATO"
$ws.Range("D22").Value = 78
$ws.Range("E22").Value = "LP1912"
$ws.Range("A23").Value = "04:18:53"
$ws.Range("B23").Value = "06:13"
$ws.Range("C23").Value = "225_HARAS DEL SUR"
$ws.Range("D23").Value = 115
$ws.Range("E23").Value = "LP1912"
$ws.Range("A24").Value = "04:18:53"
$ws.Range("B24").Value = "06:14"
$ws.Range("C24").Value = "225_HARAS DEL SUR"
$ws.Range("D24").Value = 94
$ws.Range("E24").Value = "LP1912"
$ws.Range("A25").Value = "04:40:33"
$ws.Range("B25").Value = "06:21"
$ws.Range("C25").Value = "26_HERNANDEZ"
$ws.Range("D25").Value = 101
$ws.Range("E25").Value = "LP1912"
$ws.Range("A26").Value = "04:40:33"
$ws.Range("B26").Value = "06:27"
$ws.Range("C26").Value = "23_HERNANDEZ"
$ws.Range("D26").Value = 107
$ws.Range("E26").Value = "LP1912"
$ws.Range("A27").Value = "04:40:33"
$ws.Range("B27").Value = "06:29"
$ws.Range("C27").Value = "86_EST CHICA-ESC AGRARIA"
$ws.Range("D27").Value = 109
$ws.Range("E27").Value = "LP1912"
$ws.Range("A28").Value = "04:54:03"
$ws.Range("B28").Value = "06:30"
$ws.Range("C28").Value = "86_EST CHICA-ESC AGRARIA"
$ws.Range("D28").Value = 96
$ws.Range("E28").Value = "LP1912"
$ws.Range("A29").Value = "04:40:33"
$ws.Range("B29").Value = "06:31"
$ws.Range("C29").Value = "16_SANTA ANA"
$ws.Range("D29").Value = 111
$ws.Range("E29").Value = "LP1912"
$ws.Range("A30").Value = "04:54:03"
$ws.Range("B30").Value = "06:44"
$ws.Range("C30").Value = "225_C ROCA-H SUR"
$ws.Range("D30").Value = 110
$ws.Range("E30").Value = "LP1912"
$ws.Range("A31").Value = "04:54:03"
$ws.Range("B31").Value = "06:47"
$ws.Range("C31").Value = "215C_EL PATO"
$ws.Range("D31").Value = 113
$ws.Range("E31").Value = "LP1912"

### Sheet: LP1912-215
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 04:54:03"
$ws.Range("A3").Value = "Total filas: 7"
$ws.Range("A6").Value = "04:18:53"
$ws.Range("B6").Value = "04:45"
$ws.Range("C6").Value = "215A_EL PATO"
$ws.Range("D6").Value = 27
$ws.Range("E6").Value = "LP1912"
$ws.Range("A7").Value = "04:40:33"
$ws.Range("B7").Value = "04:46"
$ws.Range("C7").Value = "215A_EL PATO"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = "LP1912"
$ws.Range("A8").Value = "04:18:53"
$ws.Range("B8").Value = "05:34"
$ws.Range("C8").Value = "215B_EL PATO"
$ws.Range("D8").Value = 76
$ws.Range("E8").Value = "LP1912"
$ws.Range("A9").Value = "04:54:03"
$ws.Range("B9").Value = "05:35"
$ws.Range("C9").Value = "215B_EL PATO"
$ws.Range("D9").Value = 41
$ws.Range("E9").Value = "LP1912"
$ws.Range("A10").Value = "04:18:53"
$ws.Range("B10").Value = "06:11"
$ws.Range("C10").Value = "215A_EL PATO"
$ws.Range("D10").Value = 113
$ws.Range("E10").Value = "LP1912"
$ws.Range("A11").Value = "04:54:03"
$ws.Range("B11").Value = "06:12"
$ws.Range("C11").Value = "215A_EL PATO"
$ws.Range("D11").Value = 78
$ws.Range("E11").Value = "LP1912"
$ws.Range("A12").Value = "04:54:03"
$ws.Range("B12").Value = "06:47"
$ws.Range("C12").Value = "215C_EL PATO"
$ws.Range("D12").Value = 113
$ws.Range("E12").Value = "LP1912"

### Sheet: 6203-6173
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 04:54:03"
$ws.Range("A3").Value = "Total filas: 5"
$ws.Range("A6").Value = "04:18:53"
$ws.Range("B6").Value = "05:43"
$ws.Range("C6").Value = "215A_LA PLATA"
$ws.Range("D6").Value = 85
$ws.Range("E6").Value = "L6173"
$ws.Range("A7").Value = "04:40:33"
$ws.Range("B7").Value = "05:44"
$ws.Range("C7").Value = "215A_LA PLATA"
$ws.Range("D7").Value = 64
$ws.Range("E7").Value = "L6173"
$ws.Range("A8").Value = "04:18:53"
$ws.Range("B8").Value = "06:08"
$ws.Range("C8").Value = "215A_LA PLATA"
$ws.Range("D8").Value = 110
$ws.Range("E8").Value = "L6173"
$ws.Range("A9").Value = "04:40:33"
$ws.Range("B9").Value = "06:09"
$ws.Range("C9").Value = "215A_LA PLATA"
$ws.Range("D9").Value = 89
$ws.Range("E9").Value = "L6173"
$ws.Range("A10").Value = "04:40:33"
$ws.Range("B10").Value = "06:33"
$ws.Range("C10").Value = "215C_LA PLATA"
$ws.Range("D10").Value = 113
$ws.Range("E10").Value = "L6203"

# trailing no-op so the script doesn't end on an array-valued Range assignment
Write-Host "Horarios actualizados Línea 141 - 253"
